{"js": "// The template paragraph holds a single Word field whose code is\n//  \" m:''.emptyList() \" (an M2Doc script field). The edit rewrites that\n// field away into plain literal text runs that spell out the same\n// script, delimited with \"{\" / \"}\" instead of field begin/end chars:\n//   <w:r><w:t>{</w:t></w:r>\n//   <w:r><w:t>m</w:t></w:r>\n//   <w:r><w:t>:''.emptyList()</w:t></w:r>\n//   <w:r><w:t xml:space=\"preserve\">}</w:t></w:r>\n\nconst fields = context.document.body.fields;\nfields.load(\"items\");\nawait context.sync();\n\nif (fields.items.length === 0) {\n  throw new Error(\"Expected a field containing the M2Doc script code.\");\n}\n\nconst field = fields.items[0];\n\n// The field result range lives inside the paragraph we need to rewrite;\n// grab that paragraph before the field (and its begin/instrText/end runs)\n// is removed.\nconst paragraph = field.result.paragraphs.getFirst();\nconst paragraphRange = paragraph.getRange();\n\n// Removing the field deletes the fldChar begin/end runs and the\n// instrText runs that made up \" m:''.emptyList() \", leaving the\n// paragraph empty.\nfield.delete();\nawait context.sync();\n\n// Re-insert the same script text as plain runs, each run holding one of\n// the former instrText chunks, but wrapped in \"{\" ... \"}\" instead of the\n// field delimiters/spaces.\nconst flatOpcXml =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body><w:p>\" +\n  \"<w:r><w:t>{</w:t></w:r>\" +\n  \"<w:r><w:t>m</w:t></w:r>\" +\n  \"<w:r><w:t>:''.emptyList()</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n  \"</w:p></w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\nparagraphRange.insertOoxml(flatOpcXml, \"Replace\");\nawait context.sync();\n", "ps1": "# The template paragraph holds a single Word field whose code is\n#  \" m:''.emptyList() \" (an M2Doc script field). The edit rewrites that\n# field away into plain literal text runs that spell out the same\n# script, delimited with \"{\" / \"}\" instead of field begin/end chars:\n#   <w:r><w:t>{</w:t></w:r>\n#   <w:r><w:t>m</w:t></w:r>\n#   <w:r><w:t>:''.emptyList()</w:t></w:r>\n#   <w:r><w:t xml:space=\"preserve\">}</w:t></w:r>\n\n$d = $word.ActiveDocument\n\n$f = $d.Fields.Item(1)\n$fStart = $f.Code.Start\n$fEnd = $f.Code.End\n\n# Locate the paragraph that contains the field (Field.Result/Field.Code's\n# own .Paragraphs collection is unreliable here, so resolve it by scanning\n# the document's paragraph ranges for the one that spans the field).\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $pr = $p.Range\n    if ($pr.Start -le $fStart -and $pr.End -ge $fEnd) {\n        $targetIndex = $i\n    }\n}\n\n# Removing the field deletes the fldChar begin/end runs and the\n# instrText runs that made up \" m:''.emptyList() \", leaving the\n# paragraph empty.\n$f.Delete()\n\n$p2 = $d.Paragraphs.Item($targetIndex)\n$r = $p2.Range\n\n# Re-insert the same script text as plain runs, each run holding one of\n# the former instrText chunks, but wrapped in \"{\" ... \"}\" instead of the\n# field delimiters/spaces.\n$flatOpc = '<?xml version=\"1.0\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:' + \"''\" + '.emptyList()</w:t></w:r><w:r><w:t xml:space=\"preserve\">}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$r.InsertXML($flatOpc)\n"}
